$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -105.5674577863067
$ws.Range("C2").Value = 0.09916957484935195
$ws.Range("D2").Value = 237.515629173

$ws.Range("B3").Value = -105.84488407880603
$ws.Range("C3").Value = 0.09920532363664455
$ws.Range("D3").Value = 192.958196244

$ws.Range("B4").Value = -105.4938270439315
$ws.Range("C4").Value = 0.09013521694606508
$ws.Range("D4").Value = 210.486877388

$ws.Range("B5").Value = -104.60548337657738
$ws.Range("C5").Value = 0.09852701461338081
$ws.Range("D5").Value = 297.926080519

$ws.Range("B6").Value = -103.49903687190475
$ws.Range("C6").Value = 0.08675010732014433
$ws.Range("D6").Value = 188.201400343

$ws.Range("B7").Value = -104.63902673475309
$ws.Range("C7").Value = 0.08112826102326275
$ws.Range("D7").Value = 217.754772853

$ws.Range("B8").Value = -102.87894895186048
$ws.Range("C8").Value = 0.09775684673713221
$ws.Range("D8").Value = 212.727060264

$ws.Range("B9").Value = -105.10124785893194
$ws.Range("C9").Value = 0.09898214076070072
$ws.Range("D9").Value = 186.591938283

$ws.Range("B10").Value = -105.9960833017315
$ws.Range("C10").Value = 0.07092340862122525
$ws.Range("D10").Value = 221.48002542

$ws.Range("B11").Value = -102.78699673530178
$ws.Range("C11").Value = 0.09590964376430258
$ws.Range("D11").Value = 183.712251881

